$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'29.760.30"
$ws.Cells.Item(2, 5).Value = "'  -0.76%  "

$ws.Cells.Item(3, 4).Value = "'1.889.09"
$ws.Cells.Item(3, 5).Value = "'  -0.97%  "

$ws.Cells.Item(4, 5).Value = "'  +0.12%  "

$ws.Cells.Item(5, 4).Value = "'0.7921"
$ws.Cells.Item(5, 5).Value = "'  -4.76%  "

$ws.Cells.Item(6, 4).Value = "'241.07"
$ws.Cells.Item(6, 5).Value = "'  -0.30%  "

$ws.Cells.Item(7, 4).Value = "'1.000"
$ws.Cells.Item(7, 5).Value = "'  +0.05%  "

$ws.Cells.Item(8, 4).Value = "'0.3150"
$ws.Cells.Item(8, 5).Value = "'  -1.48%  "

$ws.Cells.Item(9, 4).Value = "'25.49"
$ws.Cells.Item(9, 5).Value = "'  -4.21%  "

$ws.Cells.Item(10, 4).Value = "'0.06989"
$ws.Cells.Item(10, 5).Value = "'  -0.25%  "

$ws.Cells.Item(11, 4).Value = "'0.08031"
$ws.Cells.Item(11, 5).Value = "'  +0.19%  "

$ws.Cells.Item(12, 4).Value = "'0.7565"
$ws.Cells.Item(12, 5).Value = "'  +1.00%  "

$ws.Cells.Item(13, 4).Value = "'1.904.41"
$ws.Cells.Item(13, 5).Value = "'  -0.20%  "

$ws.Cells.Item(14, 4).Value = "'5.285"
$ws.Cells.Item(14, 5).Value = "'  +1.63%  "

$ws.Cells.Item(15, 5).Value = "'  -0.67%  "

$ws.Cells.Item(16, 4).Value = "'29.776.39"
$ws.Cells.Item(16, 5).Value = "'  -0.68%  "

$ws.Cells.Item(17, 4).Value = "'13.75"
$ws.Cells.Item(17, 5).Value = "'  -2.65%  "

$ws.Cells.Item(18, 4).Value = "'5.903"
$ws.Cells.Item(18, 5).Value = "'  +0.33%  "

$ws.Cells.Item(19, 4).Value = "'243.65"
$ws.Cells.Item(19, 5).Value = "'  -0.50%  "

$ws.Cells.Item(20, 4).Value = "'0.000007654"
$ws.Cells.Item(20, 5).Value = "'  -1.62%  "

$ws.Cells.Item(21, 5).Value = "'  +0.18%  "

$ws.Cells.Item(22, 4).Value = "'2.152.94"
$ws.Cells.Item(22, 5).Value = "'  -0.45%  "

$ws.Cells.Item(23, 4).Value = "'8.098"
$ws.Cells.Item(23, 5).Value = "'  +16.34%  "

$ws.Cells.Item(24, 4).Value = "'1.001"
$ws.Cells.Item(24, 5).Value = "'  +0.10%  "

$ws.Cells.Item(25, 5).Value = "'  +1.75%  "

$ws.Cells.Item(26, 4).Value = "'9.262"
$ws.Cells.Item(26, 5).Value = "'  +0.38%  "

$ws.Cells.Item(27, 4).Value = "'163.48"
$ws.Cells.Item(27, 5).Value = "'  -3.49%  "

$ws.Cells.Item(28, 4).Value = "'18.58"
$ws.Cells.Item(28, 5).Value = "'  -1.82%  "

$ws.Cells.Item(29, 4).Value = "'2.041"
$ws.Cells.Item(29, 5).Value = "'  -2.02%  "

$ws.Cells.Item(30, 5).Value = "'  +0.90%  "

$ws.Cells.Item(31, 4).Value = "'1.532"
$ws.Cells.Item(31, 5).Value = "'  +1.34%  "

$ws.Cells.Item(32, 4).Value = "'4.372"
$ws.Cells.Item(32, 5).Value = "'  +1.75%  "

$ws.Cells.Item(33, 4).Value = "'0.05663"
$ws.Cells.Item(33, 5).Value = "'  +0.99%  "

$ws.Cells.Item(34, 4).Value = "'4.045"
$ws.Cells.Item(34, 5).Value = "'  -0.92%  "

$ws.Cells.Item(35, 4).Value = "'1.260"
$ws.Cells.Item(35, 5).Value = "'  -0.83%  "

$ws.Cells.Item(36, 4).Value = "'0.7323"
$ws.Cells.Item(36, 5).Value = "'  -0.12%  "

$ws.Cells.Item(37, 4).Value = "'0.9995"
$ws.Cells.Item(37, 5).Value = "'  +0.14%  "

$ws.Cells.Item(38, 4).Value = "'2.588"
$ws.Cells.Item(38, 5).Value = "'  -4.49%  "

$ws.Cells.Item(39, 4).Value = "'0.01898"
$ws.Cells.Item(39, 5).Value = "'  -1.34%  "

$ws.Cells.Item(40, 4).Value = "'2.776"
$ws.Cells.Item(40, 5).Value = "'  -0.43%  "

$ws.Cells.Item(41, 5).Value = "'  -1.07%  "

$ws.Cells.Item(42, 4).Value = "'72.14"
$ws.Cells.Item(42, 5).Value = "'  -0.31%  "

$ws.Cells.Item(43, 4).Value = "'5.806"
$ws.Cells.Item(43, 5).Value = "'  -2.98%  "

$ws.Cells.Item(44, 5).Value = "'  +0.17%  "

$ws.Cells.Item(45, 4).Value = "'0.8379"
$ws.Cells.Item(45, 5).Value = "'  -0.39%  "

$ws.Cells.Item(46, 4).Value = "'102.24"
$ws.Cells.Item(46, 5).Value = "'  +1.25%  "

$ws.Cells.Item(47, 4).Value = "'1.018.66"
$ws.Cells.Item(47, 5).Value = "'  +2.71%  "

$ws.Cells.Item(48, 4).Value = "'1.850"
$ws.Cells.Item(48, 5).Value = "'  -2.42%  "

$ws.Cells.Item(49, 4).Value = "'9.830"
$ws.Cells.Item(49, 5).Value = "'  +0.98%  "

$ws.Cells.Item(50, 4).Value = "'7.422"
$ws.Cells.Item(50, 5).Value = "'  -2.23%  "

$ws.Cells.Item(51, 4).Value = "'2.052.02"
$ws.Cells.Item(51, 5).Value = "'  -0.50%  "

$ws.Range("D2:E51").Style = "Normal"
